$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 2 - shape "Elements de contexte sur le batiment": shrink height,
#    add explicit text insets.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(5)
$sh2.Height = 25.112165451049805
$sh2.TextFrame.MarginLeft = 36000/12700
$sh2.TextFrame.MarginTop = 36000/12700
$sh2.TextFrame.MarginRight = 36000/12700
$sh2.TextFrame.MarginBottom = 36000/12700

# ---------------------------------------------------------------------------
# 2) Slide 3 - shape "Energie et consommations": shrink height, shrink insets.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$sh3.Height = 25.112165451049805
$sh3.TextFrame.MarginLeft = 36000/12700
$sh3.TextFrame.MarginTop = 36000/12700
$sh3.TextFrame.MarginRight = 36000/12700
$sh3.TextFrame.MarginBottom = 36000/12700

# ---------------------------------------------------------------------------
# 3) Slide 5 - shape "Descriptif de l'enveloppe thermique": shrink height,
#    shrink insets.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$sh5.Height = 25.112165451049805
$sh5.TextFrame.MarginLeft = 36000/12700
$sh5.TextFrame.MarginTop = 36000/12700
$sh5.TextFrame.MarginRight = 36000/12700
$sh5.TextFrame.MarginBottom = 36000/12700

# ---------------------------------------------------------------------------
# 4) Slide 7 - "Descriptif du chauffage" textbox relocated: it is deleted
#    from its original spot (just before the "tableauEmetteurs" graphic
#    frame) and an equivalent textbox (copied from the already-updated
#    slide 5 box, which carries matching paragraph/run formatting) is
#    pasted at the end of the shape tree with the new size/position.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$oldChauffage = $s7.Shapes.Item(2)
$oldChauffage.Delete()

$sh5.Copy()
$newChauffage = $s7.Shapes.Paste().Item(1)
$newChauffage.Name = "Descriptif du chauffage"
$newChauffage.Left = 506.8365173339844
$newChauffage.Top = 414.930908203125
$newChauffage.Width = 433.5002136230469
$newChauffage.Height = 25.112165451049805
$newChauffage.TextFrame.TextRange.Text = " `t"
